$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 6.041599999999989
$ws.Range("D5").Value = -8.461999999999993
$ws.Range("E7").Value = 12.09749999999999
$ws.Range("D9").Value = -8.663900000000002
$ws.Range("D11").Value = -8.47000000000001
$ws.Range("E11").Value = 13.30379999999999
$ws.Range("B21").Value = 5.620899999999993
$ws.Range("D21").Value = -7.498400000000003
$ws.Range("E21").Value = 13.66989999999999
$ws.Range("B23").Value = 6.506899999999997
$ws.Range("B25").Value = 5.724499999999995
